$d = $word.ActiveDocument

# --- Edit 1: Table cell — mark the "UI/UX Design Specification" row's
#     "Comments" cell (row 7, col 4) as "Drafted" ---
$t1 = $d.Tables.Item(1)
$cell = $t1.Cell(7, 4)
$cell.Range.Text = "Drafted"

# --- Edit 2: Default page header — merge the tab run and the
#     "R-SRAFVP SRS" text run into a single run (keeping <w:tab/> as a
#     distinct child so it still renders as a real tab, not literal text) ---
$sec = $d.Sections.First
$hdr = $sec.Headers.Item(1)
$p = $hdr.Range.Paragraphs.Item(1)
$r = $p.Range
$r.MoveEnd(1, -1)

$xml = @'
<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:tab/><w:t>R-SRAFVP SRS</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$r.InsertXML($xml)
